# Updates cryptos list price/volume figures to the latest scraped values.
# For cells whose new text looks like a plain number (e.g. "0.7114"),
# the cell is forced to Text format before assignment and the format is
# cleared again afterwards so the stored value stays a literal string
# (matching the source data's inlineStr text cells) instead of being
# auto-coerced into a numeric value by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.288.06'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.869.92'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7114'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3106'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07726'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08396'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = '1.883.53'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.218'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7118'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").Value = '29.300.62'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008193'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.51'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = '2.123.73'
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.863'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.55%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.81'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.008'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.511'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.401'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.287'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05173'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.917'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7735'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.170'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.689'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.713'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = '1.159.75'
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.409'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.20'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8906'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("E44").Value = '  +2.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9999'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = '2.020.52'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5190'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.792'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.385'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000120'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("E51").Value = '  -0.67%  '
